$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 53
$ws.Range("F2").Value = 35
$ws.Range("H2").Value = 47

$ws.Range("E15").Value = 167
$ws.Range("F15").Value = 92
$ws.Range("H15").Value = 133

$ws.Range("E17").Value = 127

$ws.Range("E18").Value = 119
$ws.Range("F18").Value = 53
$ws.Range("H18").Value = 89

$ws.Range("E19").Value = 65
$ws.Range("F19").Value = 35
$ws.Range("H19").Value = 48

$ws.Range("E22").Value = 5

$ws.Range("E23").Value = 5

$ws.Range("E24").Value = 24

$ws.Range("E32").Value = 23

$ws.Range("E36").Value = 111
$ws.Range("F36").Value = 51
$ws.Range("H36").Value = 83

$ws.Range("E41").Value = 43

$ws.Range("E52").Value = 9

$ws.Range("E60").Value = 20

$ws.Range("E62").Value = 50

$ws.Range("E63").Value = 39

$ws.Range("E64").Value = 36

$ws.Range("E65").Value = 34

$ws.Range("E68").Value = 18

$ws.Range("E70").Value = 47
$ws.Range("F70").Value = 23
$ws.Range("H70").Value = 35

$ws.Range("E73").Value = 30

$ws.Range("F81").Value = 13
$ws.Range("H81").Value = 18

$ws.Range("E82").Value = 17
$ws.Range("F82").Value = 3
$ws.Range("H82").Value = 9

$ws.Range("E87").Value = 17
$ws.Range("F87").Value = 6
$ws.Range("H87").Value = 13
